# Add a new data row (row 32) at the bottom of the "Relatórios" report
# sheet, mirroring the layout of the existing rows (columns C:L).
#
# Every value in C:L on this sheet is stored as text (even the
# numeric-looking Part Number / Semana / Número de Relatório columns), so
# we briefly force a "Text" number format before writing those values —
# otherwise Excel would auto-coerce strings like "53490059" or "1" into
# numbers. The number format is reset back to the default "Normal" style
# immediately after, so the new cells end up with no special formatting,
# just like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32
$rng = $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 12))
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 3).Value  = "53490059"
$ws.Cells.Item($row, 4).Value  = "coluna  do conj transversal traseiro ld"
$ws.Cells.Item($row, 5).Value  = "1"
$ws.Cells.Item($row, 6).Value  = "QUALIDADE"
$ws.Cells.Item($row, 7).Value  = "izaac"
$ws.Cells.Item($row, 8).Value  = "2º TURNO"
$ws.Cells.Item($row, 9).Value  = "METRASCAN"
$ws.Cells.Item($row, 10).Value = "INSP LAYOUT"
$ws.Cells.Item($row, 11).Value = "dsxfvgrsegt gt"
$ws.Cells.Item($row, 12).Value = "C2025.0030"

$rng.Style = "Normal"
